$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts all existing columns
# (and merged ranges) one place to the right and adds the "Match ID"
# shared string automatically once we set the value below.
$ws.Range("A1").EntireColumn.Insert()

# Header label for the new first column.
$ws.Range("A2").Value = "Match ID"
$ws.Range("A2").Font.Bold = $true

# Hidden separator row keeps the bold style with no value.
$ws.Range("A3").Font.Bold = $true

# Fill the new "Match ID" column for every visible/hidden data row
# (rows 4-18) with the match id, using the same bold style as A2/A3.
$ws.Range("A4:A18").Value = 12
$ws.Range("A4:A18").Font.Bold = $true

# The trailing summary row keeps the default (unbold) style.
$ws.Range("A19").Value = 12
# Writing into this hidden row marks it for a height recalculation; re-fit
# it so it doesn't pick up an explicit row height it never had before.
[void]$ws.Rows.Item(19).AutoFit()

# Restore the selection to the new Match ID column's data range.
[void]$ws.Range("A2:A18").Select()
